$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.770.42'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '2.045.26'
$ws.Range('E3').Value = '  +1.06%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''227.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').Value = '''0.609'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.39%  '
$ws.Range('D7').Value = '''60.29'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.08%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '''0.375'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.00%  '
$ws.Range('D10').Value = '''0.0838'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.31%  '
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('D12').Value = '2.348.34'
$ws.Range('E12').Value = '  +0.99%  '
$ws.Range('D13').Value = '''14.37'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.21%  '
$ws.Range('D14').Value = '''21.39'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.28%  '
$ws.Range('D15').Value = '''5.47'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.57%  '
$ws.Range('E16').Value = '  +1.04%  '
$ws.Range('D17').Value = '2.048.13'
$ws.Range('E17').Value = '  +1.67%  '
$ws.Range('D18').Value = '37.746.13'
$ws.Range('E18').Value = '  +0.23%  '
$ws.Range('E19').Value = '  -1.57%  '
$ws.Range('D20').Value = '''69.39'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.38%  '
$ws.Range('E21').Value = '  +0.59%  '
$ws.Range('D22').Value = '''223.01'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.67%  '
$ws.Range('E23').Value = '  +0.44%  '
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').Value = '''2.27'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.33%  '
$ws.Range('D26').Value = '''168.97'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.28%  '
$ws.Range('D27').Value = '''9.32'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.99%  '
$ws.Range('E28').Value = '  -0.32%  '
$ws.Range('E29').Value = '  -0.65%  '
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('E31').Value = '  -0.84%  '
$ws.Range('E32').Value = '  +8.80%  '
$ws.Range('E33').Value = '  -1.13%  '
$ws.Range('E34').Value = '  +0.37%  '
$ws.Range('D35').Value = '''4.49'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.23%  '
$ws.Range('E36').Value = '  +3.64%  '
$ws.Range('E37').Value = '  +4.46%  '
$ws.Range('E38').Value = '  +6.91%  '
$ws.Range('E39').Value = '  -0.09%  '
$ws.Range('D40').Value = '''18.00'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.14%  '
$ws.Range('D41').Value = '1.532.45'
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('D42').Value = '''97.85'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.27%  '
$ws.Range('E43').Value = '  -0.67%  '
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('E45').Value = '  +7.49%  '
$ws.Range('D46').Value = '''0.0897'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.20%  '
$ws.Range('E48').Value = '  +0.33%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = '''7.09'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.63%  '
$ws.Range('B50').Value = 'MXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D50').Value = '''2.94'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.38%  '
$ws.Range('D51').Value = '2.236.79'
$ws.Range('E51').Value = '  +0.95%  '
